$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.115.33'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.40%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.223.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.15%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.86%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.621'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.27'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '

# Row 8
$ws.Range('E8').Value = '  +0.08%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.617'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.30%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.58%  '

# Row 11
$ws.Range('E11').Value = '  +2.64%  '

# Row 12
$ws.Range('E12').Value = '  -0.01%  '

# Row 13
$ws.Range('E13').Value = '  +0.59%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.553.35'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.26%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.850'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.75%  '

# Row 16
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.59%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.218.48'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.61%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.946.24'
$ws.Range('D18').Style = 'Normal'

# Row 19
$ws.Range('E19').Value = '  +11.95%  '

# Row 20
$ws.Range('E20').Value = '  +0.84%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.56'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.84%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.80'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +23.45%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.09%  '

# Row 24
$ws.Range('E24').Value = '  -6.39%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.74%  '

# Row 26
$ws.Range('E26').Value = '  +0.02%  '

# Row 27
$ws.Range('E27').Value = '  -1.19%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.82%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.15'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.56%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.77'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.68%  '

# Row 31
$ws.Range('E31').Value = '  +1.56%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.70'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +11.07%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0802'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.58%  '

# Row 34
$ws.Range('E34').Value = '  +0.12%  '

# Row 35
$ws.Range('E35').Value = '  -4.27%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '29.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.42%  '

# Row 37
$ws.Range('E37').Value = '  -4.64%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0305'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.49%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.11%  '

# Row 40
$ws.Range('B40').Value = 'MultiversX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '66.41'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.34%  '

# Row 41
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.14'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.06%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.66'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.26%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.202'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.44%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.13%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.83'
$ws.Range('D45').Style = 'Normal'

# Row 46
$ws.Range('E46').Value = '  +0.28%  '

# Row 47
$ws.Range('E47').Value = '  +6.16%  '

# Row 48
$ws.Range('E48').Value = '  +0.32%  '

# Row 49
$ws.Range('E49').Value = '  +0.40%  '

# Row 50
$ws.Range('E50').Value = '  +0.15%  '

# Row 51
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.428.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.24%  '
